$wb = $excel.ActiveWorkbook

# Sheet "展览" (Sheet1) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 181
$ws1.Range("F5").Value = 5001
$ws1.Range("F8").Value = 9
$ws1.Range("F9").Value = 553
$ws1.Range("F10").Value = 513
$ws1.Range("F13").Value = 1392
$ws1.Range("F14").Value = 3663
$ws1.Range("F16").Value = 136
$ws1.Range("F17").Value = 120
$ws1.Range("F18").Value = 82
$ws1.Range("F19").Value = 2655
$ws1.Range("F20").Value = 131
$ws1.Range("F22").Value = 86
$ws1.Range("F25").Value = 57
$ws1.Range("F26").Value = 129

# Sheet "全部类型" (Sheet4) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 181
$ws4.Range("F6").Value = 5001
$ws4.Range("F9").Value = 9
$ws4.Range("F10").Value = 553
$ws4.Range("F11").Value = 513
$ws4.Range("F14").Value = 1392
$ws4.Range("F15").Value = 3663
$ws4.Range("F17").Value = 136
$ws4.Range("F18").Value = 120
$ws4.Range("F19").Value = 82
$ws4.Range("F20").Value = 2655
$ws4.Range("F21").Value = 131
$ws4.Range("F23").Value = 86
$ws4.Range("F26").Value = 57
$ws4.Range("F27").Value = 129
